$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 16: task "BSJF240124 / 蒋镥" description + expected/actual dates
# ---------------------------------------------------------------------
# D16: "预测结合转录子" -> "结合转录因子"
$ws.Range("D16").Value = "结合转录因子"

# G16 (预计完成) / H16 (实际完成): 45741 -> 45743  (2025-03-25 -> 2025-03-27)
$ws.Range("G16").Value = 45743
$ws.Range("H16").Value = 45743

# ---------------------------------------------------------------------
# Row 17: previously a blank placeholder row (only carried cell styles,
# no content) - fill it in with a new "BS.develop" task entry.
# ---------------------------------------------------------------------

# B17 used to be a style-only blank cell. The target state is an empty
# TEXT cell (same as e.g. B3), not a numeric blank. Snapshot B17's
# pristine formatting first so we can restore it after forcing the
# text type below (a leading apostrophe is the standard Excel idiom
# for "treat this entry as text"; it also marks the cell quote-prefixed,
# which we undo by re-applying the original formatting).
$ws.Range("B17").Copy()
$ws.Range("K17").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# E17:H17 are date cells; reuse the same date-number-format already used
# by row 16's E:H columns (style 34) instead of letting a fresh
# NumberFormat assignment fabricate a brand-new (duplicate) style.
$ws.Range("E16:H16").Copy()
$ws.Range("E17:H17").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A17").Value = "BS.develop"
$ws.Range("B17").Value = "'"
$ws.Range("C17").Value = "模块开发"
$ws.Range("D17").Value = "转录因子数据获取和可视化"
$ws.Range("E17").Value = 45742
$ws.Range("F17").Value = 45743
$ws.Range("G17").Value = 45743
$ws.Range("H17").Value = 45743
$ws.Range("I17").Value = "抓取 hTFtarget 数据库程序，转录因子可视化程序"

# Restore B17's original formatting (undoing the quote-prefix style the
# apostrophe trick introduced) then remove the scratch cell.
$ws.Range("K17").Copy()
$ws.Range("B17").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("K17").Clear()
